# issue #5: add legislator_id, name, date into dataframe
#
# The "股票" (stocks) sheet - the 6th worksheet tab in this per-legislator
# workbook - gains three trailing columns (date, legislator_name,
# legislator_id) that tag every existing data row with the filing's
# metadata, the same way the other per-legislator export sheets do.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

$legislatorName = "許忠信"
$legislatorId = 1749
$reportDate = "2012-04-23"
$lastRow = 10

# --- Header row (row 1): copy the existing bold/bordered header format
# from the last header cell (G1) onto the three new header cells, then
# stamp in the new column names.
$ws.Range("G1").Copy()
$ws.Range("H1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("H1").Value = "date"
$ws.Range("I1").Value = "legislator_name"
$ws.Range("J1").Value = "legislator_id"

# The date column must stay a literal text value ("2012-04-23"), not get
# auto-converted into a date serial number, so force a text number format
# on the whole column range before writing into it.
$ws.Range("H2:H" + $lastRow).NumberFormat = "@"

# --- Data rows (2-10): every row gets the same report date, legislator
# name and legislator id.
for ($row = 2; $row -le $lastRow; $row++) {
    $ws.Range("H$row").Value = $reportDate
    $ws.Range("I$row").Value = $legislatorName
    $ws.Range("J$row").Value = $legislatorId
}
